$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last data row (row 6) - data now ends at row 5
$ws.Rows.Item(6).Delete()

# Update row 5 values to the new "custom accuracy" (2-decimal) readings
$ws.Range("B5").Value = 9.98
$ws.Range("C5").Value = 7.74
$ws.Range("D5").Value = 0.48
$ws.Range("E5").Value = 21.59
$ws.Range("F5").Value = 17.74
$ws.Range("G5").Value = 8.130000000000001
$ws.Range("H5").Value = 29.02
$ws.Range("I5").Value = 12.1
$ws.Range("J5").Value = 5.46
$ws.Range("K5").Value = 8.550000000000001
$ws.Range("L5").Value = 9.4
$ws.Range("M5").Value = 9.029999999999999
$ws.Range("N5").Value = 2.51
$ws.Range("O5").Value = 7.77
$ws.Range("P5").Value = 11.24
$ws.Range("Q5").Value = 6.51
$ws.Range("R5").Value = 0.07000000000000001
$ws.Range("S5").Value = 0.29
$ws.Range("T5").Value = 113.11
$ws.Range("U5").Value = 21.85
$ws.Range("V5").Value = 7.17
$ws.Range("W5").Value = 14.7
$ws.Range("X5").Value = 8.109999999999999
$ws.Range("Y5").Value = 1.05
$ws.Range("Z5").Value = 14.78
$ws.Range("AA5").Value = 6.29
$ws.Range("AB5").Value = 6.16
$ws.Range("AC5").Value = 6.51
$ws.Range("AD5").Value = 9.609999999999999
$ws.Range("AE5").Value = 0.07000000000000001
$ws.Range("AF5").Value = 26.42
$ws.Range("AG5").Value = 4.24
$ws.Range("AH5").Value = 8.970000000000001

# Narrow a few data columns by one unit (P, T, W) to match new column-width profile
$ws.Columns.Item(16).ColumnWidth = 6.17
$ws.Columns.Item(20).ColumnWidth = 7.17
$ws.Columns.Item(23).ColumnWidth = 6.17

